$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 323, shifting existing rows 323:361 down to 324:362
$ws.Rows("323:323").Insert()

# Populate the newly inserted row 323 with its data
$ws.Range("A323").Value = 8
$ws.Range("B323").Value = "Terminal La Palmera de La Serena"
$ws.Range("C323").Value = "Coquimbo"
$ws.Range("D323").Value = 44946
$ws.Range("E323").Value = 4
$ws.Range("F323").Value = 100112012
$ws.Range("G323").Value = "Espinaca"
$ws.Range("H323").Value = "Sin especificar"
$ws.Range("I323").Value = "Primera"
$ws.Range("J323").Value = 2200
$ws.Range("K323").Value = 500
$ws.Range("L323").Value = 600
$ws.Range("M323").Value = 550
$ws.Range("N323").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O323").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P323").Value = 1100
$ws.Range("Q323").Value = 0.5
$ws.Range("R323").Value = "Hortaliza"
